$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block: "Spreadsheet Integer parse5(String str)" (rows 74-76) ---
$ws.Range("C74").Value = "Spreadsheet Integer parse5(String str)"

$ws.Range("C75").Value = "Step"
$ws.Range("D75").Value = "Calc"

$ws.Range("C76").Value = "RETURN"
# D76 must contain the literal text "= parse(str) + 100" (same shared string as D20/D26)
# with the "stored as text" (quote-prefix) style used throughout the sheet.
$ws.Range("D20").Copy()
$ws.Range("D76").PasteSpecial(-4104)
$ws.Range("D70").Copy()
$ws.Range("D76").PasteSpecial(-4122)

# --- New block: "Spreadsheet Integer parse6(String str)" (rows 79-81) ---
$ws.Range("C79").Value = "Spreadsheet Integer parse6(String str)"

$ws.Range("C80").Value = "Step"
$ws.Range("D80").Value = "Calc"

$ws.Range("C81").Value = "RETURN"
$ws.Range("D20").Copy()
$ws.Range("D81").PasteSpecial(-4104)
$ws.Range("D70").Copy()
$ws.Range("D81").PasteSpecial(-4122)

$excel.CutCopyMode = $false

[void]$ws.Range("D79").Select()
